# thay đôi chiến lược chạy multi process. Sửa lại template báo cáo tổng hợp cơ sở
# Update the "Tháng 8" (August) row of the consolidated report sheet with the
# new run's last_edited_time and refreshed aggregate numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế tháng HỆ THỐNG")

# last_edited_time for the August pull moved later in the day.
# D6:D12 all shared the same timestamp string, so update them together.
$ws.Range("D6:D12").Value = "2024-08-03T20:14:00.000Z"

# Refreshed aggregate figures for Tháng 8 row (row 6)
$ws.Range("W6").Value = 14010000    # properties.Chi tiêu.number
$ws.Range("AA6").Value = 69090000   # properties.Lũy kế.formula.number
$ws.Range("AE6").Value = 83100000   # properties.Tổng doanh thu.formula.number
$ws.Range("AH6").Value = 68100000   # properties.Đã thanh toán.number
$ws.Range("AK6").Value = 8          # properties.Số lượng đơn.number
$ws.Range("AN6").Value = 15000000   # properties.Thu nợ.number
$ws.Range("AQ6").Value = 71100000   # properties.Đơn giá.number
